$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Force the cell to store the literal text supplied (coin-ranking numbers are
    # stored as plain strings in this sheet, e.g. "27.30" or "1.00" - letting COM
    # auto-detect the type would coerce them to numbers and drop the exact
    # formatting/trailing zeros). Resetting the style back to "Normal" afterwards
    # keeps the cell style index unchanged (style 0), matching the original file.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "63.523.70"
Set-TextValue "E2" "  +2.46%  "

# Row 3
Set-TextValue "D3" "2.473.11"
Set-TextValue "E3" "  +2.19%  "

# Row 4
Set-TextValue "E4" "  +0.15%  "

# Row 5
Set-TextValue "D5" "575.89"
Set-TextValue "E5" "  +2.30%  "

# Row 6
Set-TextValue "D6" "148.44"
Set-TextValue "E6" "  +3.63%  "

# Row 7
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.07%  "

# Row 8
Set-TextValue "E8" "  +1.76%  "

# Row 9
Set-TextValue "E9" "  +4.50%  "

# Row 10
Set-TextValue "D10" "0.154"
Set-TextValue "E10" "  +0.62%  "

# Row 11
Set-TextValue "D11" "0.363"
Set-TextValue "E11" "  +3.94%  "

# Row 12
Set-TextValue "D12" "5.34"
Set-TextValue "E12" "  +2.83%  "

# Row 13
Set-TextValue "D13" "27.30"
Set-TextValue "E13" "  +4.35%  "

# Row 14
Set-TextValue "D14" "0.0000185"
Set-TextValue "E14" "  +6.54%  "

# Row 16
Set-TextValue "D16" "63.459.71"
Set-TextValue "E16" "  +2.49%  "

# Row 17
Set-TextValue "D17" "2.492.88"
Set-TextValue "E17" "  +3.06%  "

# Row 18
Set-TextValue "D18" "11.55"
Set-TextValue "E18" "  +1.90%  "

# Row 19
Set-TextValue "D19" "7.29"
Set-TextValue "E19" "  +6.87%  "

# Row 20
Set-TextValue "B20" "BitcoinCash"
Set-TextValue "C20" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "328.98"
Set-TextValue "E20" "  +1.77%  "

# Row 21
Set-TextValue "B21" "Polkadot"
Set-TextValue "C21" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D21" "4.24"
Set-TextValue "E21" "  +2.68%  "

# Row 22
Set-TextValue "E22" "  -0.05%  "

# Row 23
Set-TextValue "E23" "  +10.70%  "

# Row 24
Set-TextValue "D24" "67.51"
Set-TextValue "E24" "  +1.15%  "

# Row 25
Set-TextValue "D25" "633.39"
Set-TextValue "E25" "  +14.34%  "

# Row 26
Set-TextValue "E26" "  +12.99%  "

# Row 27
Set-TextValue "D27" "8.75"
Set-TextValue "E27" "  -0.19%  "

# Row 28
Set-TextValue "D28" "2.595.79"
Set-TextValue "E28" "  +2.21%  "

# Row 29
Set-TextValue "D29" "1.52"
Set-TextValue "E29" "  +9.71%  "

# Row 30
Set-TextValue "D30" "8.45"
Set-TextValue "E30" "  +3.07%  "

# Row 31
Set-TextValue "E31" "  -0.22%  "

# Row 32
Set-TextValue "E32" "  -2.07%  "

# Row 33
Set-TextValue "E33" "  +1.69%  "

# Row 34
Set-TextValue "D34" "5.20"
Set-TextValue "E34" "  +9.75%  "

# Row 35
Set-TextValue "E35" "  +3.44%  "

# Row 37
Set-TextValue "D37" "0.386"
Set-TextValue "E37" "  +2.03%  "

# Row 38
Set-TextValue "D38" "5.53"
Set-TextValue "E38" "  +1.86%  "

# Row 39
Set-TextValue "D39" "18.98"
Set-TextValue "E39" "  +2.35%  "

# Row 40
Set-TextValue "D40" "1.85"
Set-TextValue "E40" "  +1.92%  "

# Row 41
Set-TextValue "D41" "147.28"
Set-TextValue "E41" "  -3.87%  "

# Row 42
Set-TextValue "D42" "2.69"
Set-TextValue "E42" "  +20.81%  "

# Row 43
Set-TextValue "E43" "  +0.50%  "

# Row 44
Set-TextValue "D44" "150.84"
Set-TextValue "E44" "  +2.39%  "

# Row 45
Set-TextValue "E45" "  +3.73%  "

# Row 46
Set-TextValue "B46" "Hedera"
Set-TextValue "C46" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D46" "0.0551"
Set-TextValue "E46" "  +4.22%  "

# Row 47
Set-TextValue "B47" "InjectiveProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D47" "21.22"
Set-TextValue "E47" "  +7.10%  "

# Row 48
Set-TextValue "D48" "0.610"
Set-TextValue "E48" "  +2.94%  "

# Row 49
Set-TextValue "E49" "  +5.38%  "

# Row 50
Set-TextValue "D50" "0.0928"
Set-TextValue "E50" "  +0.85%  "

# Row 51
Set-TextValue "D51" "0.748"
Set-TextValue "E51" "  +4.91%  "
